# Re-run of the MPI matrix-multiply benchmark results after fixing the bug
# that caused incorrect behaviour on 1 processor. The raw timings for every
# matrix-size / processor-count combination were re-measured, and the
# mislabeled "64x64" case for the "Large" matrix was corrected to "256x256".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Fix the mislabeled matrix size for the "Large" group.
$ws.Range("A6").Value = "256x256"

# Updated raw timing measurements (column C) for every run.
$ws.Range("C2").Value  = 0.006687
$ws.Range("C3").Value  = 0.010786
$ws.Range("C4").Value  = 0.011382
$ws.Range("C5").Value  = 0.155028
$ws.Range("C6").Value  = 0.136387
$ws.Range("C7").Value  = 0.151563
$ws.Range("C8").Value  = 12.181229
$ws.Range("C9").Value  = 4.399632
$ws.Range("C10").Value = 3.096313

# Re-enter the efficiency formula across the whole column in one shot so
# Excel collapses it back down into a single shared-formula group.
$ws.Range("E2:E10").Formula = "=D2/B2"

# Leave the selection where the author ended up after the edits.
$ws.Range("E14").Select()
